$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 2017.6111
$ws.Range("I41").Value = 1032.1875
$ws.Range("J41").Value = 9901
$ws.Range("K41").Value = 1032.1875
$ws.Range("L41").Value = 9901
$ws.Range("M41").Value = -592.1875
$ws.Range("N41").Value = -10781

$ws.Range("H43").Value = 5520.875
$ws.Range("I43").Value = 2386
$ws.Range("J43").Value = 6565.8335
$ws.Range("K43").Value = 2386
$ws.Range("L43").Value = 6565.8335
$ws.Range("M43").Value = -2317
$ws.Range("N43").Value = -6703.8335

$ws.Range("H61").Value = 506.33334
$ws.Range("I61").Value = 506.33334
$ws.Range("K61").Value = 1519.00002
$ws.Range("M61").Value = -1347.00002

$ws.Range("H116").Value = 10317.467
$ws.Range("I116").Value = 8896.571
$ws.Range("J116").Value = 11560.75
$ws.Range("K116").Value = 8896.571
$ws.Range("L116").Value = 11560.75
$ws.Range("M116").Value = -5454.571
$ws.Range("N116").Value = -18444.75

$ws.Range("H129").Value = 2088.7144
$ws.Range("I129").Value = 1605.8334
$ws.Range("J129").Value = 2732.5557
$ws.Range("K129").Value = 4817.5002
$ws.Range("L129").Value = 8197.667099999999
$ws.Range("M129").Value = 182.4997999999996
$ws.Range("N129").Value = -18197.6671

$ws.Range("H137").Value = 3487.6316

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 66670956
$ws.Range("I45").Value = 142858750
$ws.Range("J45").Value = 6632.625
$ws.Range("K45").Value = 142858750
$ws.Range("L45").Value = 6632.625
$ws.Range("M45").Value = -142858373
$ws.Range("N45").Value = -7386.625

$ws.Range("H74").Value = 9525102
$ws.Range("I74").Value = 11495578
$ws.Range("J74").Value = 1134
$ws.Range("K74").Value = 11495578
$ws.Range("L74").Value = 1134
$ws.Range("M74").Value = -11494704
$ws.Range("N74").Value = -2882

$ws.Range("H77").Value = 9525102
$ws.Range("I77").Value = 11495578
$ws.Range("J77").Value = 1134
$ws.Range("K77").Value = 57477890
$ws.Range("L77").Value = 5670
$ws.Range("M77").Value = -57473522
$ws.Range("N77").Value = -14406

$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws.Range("H132").Value = 2316.075
$ws.Range("I132").Value = 939.2646999999999
$ws.Range("K132").Value = 2817.7941
$ws.Range("M132").Value = -287.7941000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 8051.2104
$ws.Range("I105").Value = 1398.4286
$ws.Range("K105").Value = 1398.4286
$ws.Range("M105").Value = 348.5714

$ws.Range("H134").Value = 3006.6775
$ws.Range("I134").Value = 1488.8889
$ws.Range("J134").Value = 13251.75
$ws.Range("K134").Value = 4466.6667
$ws.Range("L134").Value = 39755.25
$ws.Range("M134").Value = -1931.6667
$ws.Range("N134").Value = -44825.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 44188.25
$ws.Range("I31").Value = 3654
$ws.Range("J31").Value = 79317.92999999999
$ws.Range("K31").Value = 3654
$ws.Range("L31").Value = 79317.92999999999
$ws.Range("M31").Value = -3359
$ws.Range("N31").Value = -79907.92999999999

$ws.Range("H34").Value = 44188.25
$ws.Range("I34").Value = 3654
$ws.Range("J34").Value = 79317.92999999999
$ws.Range("K34").Value = 3654
$ws.Range("L34").Value = 79317.92999999999
$ws.Range("M34").Value = -3452
$ws.Range("N34").Value = -79721.92999999999

$ws.Range("H62").Value = 8250.286
$ws.Range("I62").Value = 3931.6667
$ws.Range("J62").Value = 11489.25
$ws.Range("K62").Value = 3931.6667
$ws.Range("L62").Value = 11489.25
$ws.Range("M62").Value = -3307.6667
$ws.Range("N62").Value = -12737.25

$ws.Range("H65").Value = 8250.286
$ws.Range("I65").Value = 3931.6667
$ws.Range("J65").Value = 11489.25
$ws.Range("K65").Value = 19658.3335
$ws.Range("L65").Value = 57446.25
$ws.Range("M65").Value = -16538.3335
$ws.Range("N65").Value = -63686.25

$ws.Range("H103").Value = 36666.668
$ws.Range("I103").Value = 36666.668
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 36666.668
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -35494.668
$ws.Range("N103").ClearContents()

$ws.Range("H132").Value = 5226.828
$ws.Range("J132").Value = 5751.067
$ws.Range("L132").Value = 17253.201
$ws.Range("N132").Value = -22313.201

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 160287.83
$ws.Range("J37").Value = 160287.83
$ws.Range("L37").Value = 480863.49
$ws.Range("N37").Value = -481087.49

$ws.Range("H41").Value = 1003
$ws.Range("J41").Value = 1003
$ws.Range("L41").Value = 3009
$ws.Range("N41").Value = -3685

$ws.Range("H44").Value = 971.9
$ws.Range("I44").Value = 495.16666
$ws.Range("J44").Value = 1687
$ws.Range("K44").Value = 1485.49998
$ws.Range("L44").Value = 5061
$ws.Range("M44").Value = -1087.49998
$ws.Range("N44").Value = -5857

$ws.Range("H104").Value = 9000
$ws.Range("I104").Value = 8000
$ws.Range("K104").Value = 24000
$ws.Range("M104").Value = -21379

$ws.Range("H114").Value = 1312.6666
$ws.Range("I114").Value = 916.7778
$ws.Range("K114").Value = 2750.3334
$ws.Range("M114").Value = 503.6666

$ws.Range("H136").Value = 2598.2856
$ws.Range("I136").Value = 1677.8
$ws.Range("J136").Value = 4899.5
$ws.Range("K136").Value = 5033.4
$ws.Range("L136").Value = 14698.5
$ws.Range("M136").Value = 66.60000000000036
$ws.Range("N136").Value = -24898.5

$ws.Range("H140").Value = 3337.6191
$ws.Range("I140").Value = 2193.0715
$ws.Range("K140").Value = 6579.2145
$ws.Range("M140").Value = -1399.2145

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5863.875
$ws.Range("J126").Value = 12838
$ws.Range("L126").Value = 38514
$ws.Range("N126").Value = -43454

$ws.Range("H132").Value = 43258.035
$ws.Range("J132").Value = 9365.5
$ws.Range("L132").Value = 28096.5
$ws.Range("N132").Value = -33156.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3376.6667
$ws.Range("I16").Value = 2596.1538
$ws.Range("K16").Value = 2596.1538
$ws.Range("M16").Value = -2426.1538

$ws.Range("H22").Value = 18800.4
$ws.Range("J22").Value = 18800.4
$ws.Range("L22").Value = 18800.4
$ws.Range("N22").Value = -19390.4

$ws.Range("H27").Value = 18800.4
$ws.Range("J27").Value = 18800.4
$ws.Range("L27").Value = 18800.4
$ws.Range("N27").Value = -19014.4

$ws.Range("H46").Value = 5343.143
$ws.Range("J46").Value = 6000.4
$ws.Range("L46").Value = 6000.4
$ws.Range("N46").Value = -6376.4

$ws.Range("H55").Value = 1786221.5
$ws.Range("I55").Value = 2381385.2
$ws.Range("J55").Value = 730.4286
$ws.Range("K55").Value = 2381385.2
$ws.Range("L55").Value = 730.4286
$ws.Range("M55").Value = -2381212.2
$ws.Range("N55").Value = -1076.4286

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 74950
$ws.Range("J46").Value = 74950
$ws.Range("L46").Value = 74950
$ws.Range("N46").Value = -75412

$ws.Range("H100").Value = 1899.1666
$ws.Range("I100").Value = 1448.75
$ws.Range("K100").Value = 2897.5
$ws.Range("M100").Value = -2356.5

$ws.Range("H107").Value = 395.55554
$ws.Range("I107").Value = 268.66666
$ws.Range("K107").Value = 805.9999799999999
$ws.Range("M107").Value = 1114.00002

$ws.Range("H113").Value = 878.4
$ws.Range("I113").Value = 878.4
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2635.2
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -465.1999999999998
$ws.Range("N113").ClearContents()

$ws.Range("H126").Value = 5515.8945
$ws.Range("I126").Value = 5247.4707
$ws.Range("K126").Value = 15742.4121
$ws.Range("M126").Value = -13272.4121

$ws.Range("H134").Value = 74950
$ws.Range("J134").Value = 74950
$ws.Range("L134").Value = 224850
$ws.Range("N134").Value = -229920
